# Fruta / hortaliza, semanal
# Updates the weekly Espárragos (asparagus) price records for Mapocho
# Venta Directa de Santiago: dates, variety/quality, volumes, prices and
# commercialisation units for rows 2-23 are refreshed to reflect the
# latest weekly data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44160
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 9500
$ws.Range("M2").Value = 9275
$ws.Range("N2").Value = '$/bandeja 8 kilos'
$ws.Range("P2").Value = 1159
$ws.Range("Q2").Value = 8
$ws.Range("D3").Value = 44160
$ws.Range("J3").Value = 440
$ws.Range("K3").Value = 7500
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 7784
$ws.Range("N3").Value = '$/bandeja 8 kilos'
$ws.Range("P3").Value = 973
$ws.Range("Q3").Value = 8
$ws.Range("D4").Value = 44160
$ws.Range("J4").Value = 305
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 6500
$ws.Range("M4").Value = 6270
$ws.Range("N4").Value = '$/bandeja 8 kilos'
$ws.Range("P4").Value = 784
$ws.Range("Q4").Value = 8
$ws.Range("D5").Value = 44162
$ws.Range("H5").Value = 'Verde'
$ws.Range("I5").Value = 'Banquete'
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 1000
$ws.Range("D6").Value = 44162
$ws.Range("H6").Value = 'Verde'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 220
$ws.Range("K6").Value = 8500
$ws.Range("L6").Value = 8500
$ws.Range("M6").Value = 8500
$ws.Range("O6").Value = 'Región Metropolitana'
$ws.Range("P6").Value = 850
$ws.Range("D7").Value = 44162
$ws.Range("H7").Value = 'Verde'
$ws.Range("I7").Value = 'Segunda'
$ws.Range("J7").Value = 260
$ws.Range("K7").Value = 7500
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = 7500
$ws.Range("N7").Value = '$/caja 10 kilos'
$ws.Range("P7").Value = 750
$ws.Range("D8").Value = 44169
$ws.Range("H8").Value = 'Verde'
$ws.Range("I8").Value = 'Banquete'
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("P8").Value = 1000
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 95
$ws.Range("K9").Value = 7500
$ws.Range("L9").Value = 7500
$ws.Range("M9").Value = 7500
$ws.Range("P9").Value = 750
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 110
$ws.Range("K10").Value = 6500
$ws.Range("L10").Value = 6500
$ws.Range("M10").Value = 6500
$ws.Range("P10").Value = 650
$ws.Range("D11").Value = 44159
$ws.Range("I11").Value = 'Banquete'
$ws.Range("J11").Value = 180
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("O11").Value = 'Provincia de Linares'
$ws.Range("P11").Value = 1000
$ws.Range("D12").Value = 44159
$ws.Range("I12").Value = 'Primera'
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 9000
$ws.Range("N12").Value = '$/bandeja 10 kilos'
$ws.Range("P12").Value = 900
$ws.Range("D13").Value = 44159
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 320
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 7000
$ws.Range("N13").Value = '$/bandeja 10 kilos'
$ws.Range("P13").Value = 700
$ws.Range("D14").Value = 44168
$ws.Range("I14").Value = 'Banquete'
$ws.Range("J14").Value = 105
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9429
$ws.Range("P14").Value = 943
$ws.Range("D15").Value = 44168
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 290
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7500
$ws.Range("M15").Value = 7241
$ws.Range("N15").Value = '$/caja 10 kilos'
$ws.Range("P15").Value = 724
$ws.Range("D16").Value = 44168
$ws.Range("I16").Value = 'Segunda'
$ws.Range("J16").Value = 360
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 6500
$ws.Range("M16").Value = 6278
$ws.Range("N16").Value = '$/caja 10 kilos'
$ws.Range("P16").Value = 628
$ws.Range("D17").Value = 44176
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 12000
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 1200
$ws.Range("D18").Value = 44176
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Segunda'
$ws.Range("J18").Value = 100
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 9000
$ws.Range("N18").Value = '$/bandeja 10 kilos'
$ws.Range("P18").Value = 900
$ws.Range("Q18").Value = 10
$ws.Range("D19").Value = 44161
$ws.Range("I19").Value = 'Banquete'
$ws.Range("J19").Value = 260
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 9500
$ws.Range("M19").Value = 9269
$ws.Range("N19").Value = '$/caja 10 kilos'
$ws.Range("O19").Value = 'Provincia de Linares'
$ws.Range("P19").Value = 927
$ws.Range("Q19").Value = 10
$ws.Range("D20").Value = 44161
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 450
$ws.Range("K20").Value = 7500
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 7778
$ws.Range("N20").Value = '$/caja 10 kilos'
$ws.Range("O20").Value = 'Provincia de Linares'
$ws.Range("P20").Value = 778
$ws.Range("Q20").Value = 10
$ws.Range("D21").Value = 44161
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 6000
$ws.Range("L21").Value = 6500
$ws.Range("M21").Value = 6300
$ws.Range("P21").Value = 630
$ws.Range("D22").Value = 44175
$ws.Range("H22").Value = 'Sin especificar'
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 12000
$ws.Range("P22").Value = 1200
$ws.Range("D23").Value = 44175
$ws.Range("H23").Value = 'Sin especificar'
$ws.Range("J23").Value = 140
$ws.Range("K23").Value = 9000
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = 9000
$ws.Range("P23").Value = 900
